# Weekly price-list update: a new daily observation is inserted at row 56
# of the "Achicoria" sheet (pushing the existing rows 56-164 down to
# 57-165), and the newly inserted row is seeded with a copy of the row
# that lands on top of it, then the two changed fields (Fecha / Volumen)
# are overwritten with their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 56; everything below shifts down one row.
$ws.Rows("56").Insert()

# Populate the newly-blank row 56 with the same data as the row now
# sitting below it (old row 56, now at row 57), then fix up the two
# cells that actually differ for this new record.
$ws.Range("A57:R57").Copy()
$ws.Range("A56:R56").PasteSpecial()

$ws.Range("D56").Value = 45259
$ws.Range("J56").Value = 45
